$wb = $excel.ActiveWorkbook

# --- Update the formulas in m2070_overlap_fast_vcl (rows 21-32, cols R:U) ---
# Old formulas compared same-sheet row N against m2070_nooverlap_vcl row N
# (a relative "speed-up vs baseline" check). New formulas instead divide the
# same-sheet row N by the OTHER sheet's row N+20 (comparing against the newly
# added strong-scaling data block in rows 23-34 of m2070_nooverlap_vcl).
$ws2 = $wb.Worksheets.Item("m2070_overlap_fast_vcl")

$cols = @("D", "E", "F", "G")
$dstCols = @("R", "S", "T", "U")

for ($i = 0; $i -lt 12; $i++) {
    $destRow = 21 + $i
    $srcRow = 3 + $i
    $otherRow = 23 + $i
    for ($c = 0; $c -lt 4; $c++) {
        $col = $cols[$c]
        $dstCol = $dstCols[$c]
        $cell = $ws2.Range("$dstCol$destRow")
        $cell.Formula = "=$col$srcRow/m2070_nooverlap_vcl!$col$otherRow"
    }
}

# --- View-state changes ---
$ws1 = $wb.Worksheets.Item("m2070_nooverlap_vcl")
$ws1.Range("F26").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 2
$ws2.Range("R21:U32").Select()

# --- Workbook window width ---
$excel.ActiveWindow.Width = 27320
